$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.845.14'
$ws.Range("E2").Value = '  -0.92%  '
$ws.Range("D3").Value = '1.662.84'
$ws.Range("E3").Value = '  +0.30%  '
$ws.Range("E4").Value = '  -0.15%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '215.32'
$c.ClearFormats()
$ws.Range("E5").Value = '  -0.06%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.535'
$c.ClearFormats()
$ws.Range("E6").Value = '  +5.42%  '
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("E8").Value = '  +0.57%  '
$ws.Range("E9").Value = '  +0.83%  '
$ws.Range("E10").Value = '  +3.38%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0897'
$c.ClearFormats()
$ws.Range("E11").Value = '  +3.83%  '
$ws.Range("D12").Value = '1.896.98'
$ws.Range("E12").Value = '  +0.21%  '
$ws.Range("D13").Value = '1.662.84'
$ws.Range("E13").Value = '  +0.31%  '
$ws.Range("E14").Value = '  +0.07%  '
$ws.Range("E15").Value = '  +0.97%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '65.99'
$c.ClearFormats()
$ws.Range("E16").Value = '  +1.80%  '
$ws.Range("D17").Value = '26.841.91'
$ws.Range("E17").Value = '  -0.91%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '232.66'
$c.ClearFormats()
$ws.Range("E18").Value = '  -2.27%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '7.92'
$c.ClearFormats()
$ws.Range("E19").Value = '  +0.94%  '
$ws.Range("D20").Value = '0.0₃0731'
$ws.Range("E20").Value = '  +0.34%  '
$ws.Range("E21").Value = '  -0.11%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '4.40'
$c.ClearFormats()
$ws.Range("E22").Value = '  -0.40%  '
$ws.Range("E23").Value = '  -2.24%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '9.15'
$c.ClearFormats()
$ws.Range("E24").Value = '  -1.31%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '145.86'
$c.ClearFormats()
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("E26").Value = '  -0.67%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.115'
$c.ClearFormats()
$ws.Range("E27").Value = '  +1.50%  '
$ws.Range("E28").Value = '  +0.48%  '
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("E30").Value = '  -0.24%  '
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("E32").Value = '  +1.93%  '
$ws.Range("D33").Value = '1.459.25'
$ws.Range("E33").Value = '  -5.20%  '
$ws.Range("E34").Value = '  +3.60%  '
$ws.Range("E35").Value = '  +3.24%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '2.41'
$c.ClearFormats()
$ws.Range("E36").Value = '  -0.34%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.575'
$c.ClearFormats()
$ws.Range("E37").Value = '  +0.23%  '
$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.898'
$c.ClearFormats()
$ws.Range("E38").Value = '  +1.32%  '
$ws.Range("E39").Value = '  -0.07%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '5.75'
$c.ClearFormats()
$ws.Range("E40").Value = '  -3.40%  '
$ws.Range("E41").Value = '  -0.16%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '2.25'
$c.ClearFormats()
$ws.Range("E42").Value = '  -0.73%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.976'
$c.ClearFormats()
$ws.Range("E43").Value = '  +5.85%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '65.71'
$c.ClearFormats()
$ws.Range("E44").Value = '  -0.79%  '
$ws.Range("D45").Value = '1.807.74'
$ws.Range("E45").Value = '  +0.42%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.777'
$c.ClearFormats()
$ws.Range("E46").Value = '  +0.53%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '90.28'
$c.ClearFormats()
$ws.Range("E47").Value = '  +0.31%  '
$ws.Range("E48").Value = '  +0.36%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.102'
$c.ClearFormats()
$ws.Range("E49").Value = '  +4.47%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.0506'
$c.ClearFormats()
$ws.Range("E50").Value = '  +0.36%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '7.53'
$c.ClearFormats()
$ws.Range("E51").Value = '  +0.39%  '
